$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Increment the "Förändrad" (column C) date by one day for every data row (2-483).
for ($r = 2; $r -le 483; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value2 = 45182
    }
}
